$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the neighboring header cell (G1) onto the new H1 header
# cell so it reuses the same bold/bordered/centered style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column's values for each data row
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
